$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original formatting/style while forcing text (string) values,
# since source cells are inline strings and must remain text after the update
# (otherwise Excel would auto-convert numeric-looking text into real numbers).
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '51.840.29'
$ws.Range("E2").Value = '  +0.78%  '

# Row 3
$ws.Range("D3").Value = '2.813.30'
$ws.Range("E3").Value = '  +2.20%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").Value = '353.94'
$ws.Range("E5").Value = '  +6.48%  '

# Row 6
$ws.Range("E6").Value = '  -2.26%  '

# Row 7
$ws.Range("E7").Value = '  +2.36%  '

# Row 8
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  +4.01%  '

# Row 10
$ws.Range("D10").Value = '41.56'
$ws.Range("E10").Value = '  -0.10%  '

# Row 11
$ws.Range("E11").Value = '  -0.61%  '

# Row 12
$ws.Range("D12").Value = '19.96'
$ws.Range("E12").Value = '  -1.31%  '

# Row 13
$ws.Range("E13").Value = '  +1.41%  '

# Row 14
$ws.Range("E14").Value = '  +0.88%  '

# Row 15
$ws.Range("D15").Value = '3.233.53'
$ws.Range("E15").Value = '  +1.81%  '

# Row 16
$ws.Range("D16").Value = '2.812.52'
$ws.Range("E16").Value = '  +1.39%  '

# Row 17
$ws.Range("D17").Value = '0.895'
$ws.Range("E17").Value = '  +0.98%  '

# Row 18
$ws.Range("D18").Value = '51.748.56'
$ws.Range("E18").Value = '  +0.73%  '

# Row 19
$ws.Range("D19").Value = '7.38'
$ws.Range("E19").Value = '  +7.73%  '

# Row 20
$ws.Range("E20").Value = '  -1.43%  '

# Row 21
$ws.Range("D21").Value = '13.51'
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0991'
$ws.Range("E22").Value = '  +1.69%  '

# Row 23
$ws.Range("D23").Value = '269.98'
$ws.Range("E23").Value = '  -3.02%  '

# Row 24
$ws.Range("E24").Value = '  +0.24%  '

# Row 25
$ws.Range("E25").Value = '  +4.67%  '

# Row 26
$ws.Range("D26").Value = '26.69'
$ws.Range("E26").Value = '  -0.40%  '

# Row 27
$ws.Range("E27").Value = '  +0.07%  '

# Row 28
$ws.Range("E28").Value = '  +1.03%  '

# Row 29
$ws.Range("E29").Value = '  +1.71%  '

# Row 30
$ws.Range("E30").Value = '  -0.38%  '

# Row 31
$ws.Range("D31").Value = '50.61'
$ws.Range("E31").Value = '  +1.45%  '

# Row 32
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '33.74'
$ws.Range("E32").Value = '  -4.29%  '

# Row 33
$ws.Range("B33").Value = 'VeChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D33").Value = '0.0451'
$ws.Range("E33").Value = '  +28.60%  '

# Row 34
$ws.Range("E34").Value = '  +4.76%  '

# Row 35
$ws.Range("E35").Value = '  +0.82%  '

# Row 36
$ws.Range("E36").Value = '  +0.09%  '

# Row 37
$ws.Range("E37").Value = '  -0.57%  '

# Row 38
$ws.Range("D38").Value = '4.89'
$ws.Range("E38").Value = '  -1.63%  '

# Row 39
$ws.Range("D39").Value = '3.20'
$ws.Range("E39").Value = '  -1.55%  '

# Row 40
$ws.Range("D40").Value = '18.28'
$ws.Range("E40").Value = '  -4.32%  '

# Row 41
$ws.Range("D41").Value = '23.67'
$ws.Range("E41").Value = '  +1.55%  '

# Row 42
$ws.Range("E42").Value = '  +4.62%  '

# Row 43
$ws.Range("D43").Value = '126.53'
$ws.Range("E43").Value = '  -0.16%  '

# Row 44
$ws.Range("E44").Value = '  +1.25%  '

# Row 45
$ws.Range("E45").Value = '  +0.37%  '

# Row 46
$ws.Range("D46").Value = '2.079.03'
$ws.Range("E46").Value = '  -0.51%  '

# Row 47
$ws.Range("E47").Value = '  +0.22%  '

# Row 49
$ws.Range("D49").Value = '5.67'
$ws.Range("E49").Value = '  +2.68%  '

# Row 50
$ws.Range("D50").Value = '0.940'
$ws.Range("E50").Value = '  +8.58%  '

# Row 51
$ws.Range("D51").Value = '60.43'
$ws.Range("E51").Value = '  +0.82%  '

# Restore original style/number formatting for the whole data range
$dataRange.Style = $origStyle
